$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 2

$ws.Range("B3").Value = 0
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 1

$ws.Range("C4").Value = 7
$ws.Range("F4").Value = 1

$ws.Range("C5").Value = 3

$ws.Range("D6").Value = 0
$ws.Range("G6").Value = 0

$ws.Range("F7").Value = 2
